$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing
# "Opt Portfolio" column (B -> C) and "Opt Portfolio with View" column
# (C -> D) one position to the right, making room for the new
# "Initial Weights" column.
$ws.Columns("B").Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Initial Weights"

# Copy the header formatting (bold font, border, centered/top aligned)
# from the neighbouring header cell so the new header matches the look
# of the existing ones.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The column insert causes the new data cells (B2:B8) to inherit the
# bordered/bold header style from column A; clear that back to the
# default (unformatted) style used by the other data columns.
$ws.Range("B2:B8").ClearFormats()

# Populate the new "Initial Weights" values.
$ws.Range("B2").Value = 0.1
$ws.Range("B3").Value = 0.05
$ws.Range("B4").Value = 0.1
$ws.Range("B5").Value = 0.1
$ws.Range("B6").Value = 0.15
$ws.Range("B7").Value = 0.2
$ws.Range("B8").Value = 0.3
